$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for data rows 2-27 was updated from 2023-12-22 (45282)
# to 2023-12-23 (45283) by the automatic update process.
$ws.Range("C2:C27").Value = 45283
